$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the three new paragraphs at the very start of the body, before
#    the existing "2.1 - Criando o projeto..." paragraph.
# ---------------------------------------------------------------------------
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParasXml = @"
<w:p $ns>
  <w:pPr>
    <w:rPr>
      <w:highlight w:val="yellow"/>
      <w:lang w:val="pt-PT"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve">Contraints = </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
      <w:lang w:val="pt-PT"/>
    </w:rPr>
    <w:t>restriçõe</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
      <w:lang w:val="pt-PT"/>
    </w:rPr>
    <w:t>s</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:rPr>
      <w:highlight w:val="yellow"/>
      <w:lang w:val="pt-PT"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:rPr>
      <w:highlight w:val="yellow"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@

$startRange = $d.Range(0, 0)
$startRange.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# 2) Add the "HTML Preformatted" paragraph style and its linked "HTML
#    pré-formatado Caráter" character style to the style sheet.
# ---------------------------------------------------------------------------
$stP = $d.Styles.Add("HTMLpr-formatado", 1)            # wdStyleTypeParagraph
$stP.NameLocal = "HTML Preformatted"
$stP.BaseStyle = "Normal"
$stP.Priority = 99
$stP.UnhideWhenUsed = $true
$stP.ParagraphFormat.SpaceAfter = 0
$stP.ParagraphFormat.LineSpacingRule = 0                # wdLineSpaceSingle
$stP.Font.Name = "Consolas"
$stP.Font.Size = 10
$stP.Font.SizeBi = 10

$stC = $d.Styles.Add("HTMLpr-formatadoCarter", 2)       # wdStyleTypeCharacter
$stC.NameLocal = "HTML pré-formatado Caráter"
$stC.BaseStyle = "Tipodeletrapredefinidodopargrafo"
$stC.Priority = 99
$stC.Font.Name = "Consolas"
$stC.Font.Size = 10
$stC.Font.SizeBi = 10

$stP.LinkStyle = "HTMLpr-formatadoCarter"
$stC.LinkStyle = "HTMLpr-formatado"
